# Auto-generated edit script: updates cached market-data values (columns H-N)
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect
# refreshed pricing figures from the scheduled data-refresh run.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 8596.875
$ws.Range("I51").Value = 8629.333000000001
$ws.Range("J51").Value = 8589.385
$ws.Range("K51").Value = 8629.333000000001
$ws.Range("L51").Value = 8589.385
$ws.Range("M51").Value = -8145.333000000001
$ws.Range("N51").Value = -9557.385
# Row 57
$ws.Range("H57").Value = 59946.5
$ws.Range("J57").Value = 59946.5
$ws.Range("L57").Value = 179839.5
$ws.Range("N57").Value = -180837.5
# Row 62
$ws.Range("H62").Value = 2732.7778
$ws.Range("I62").Value = 2248
$ws.Range("K62").Value = 2248
$ws.Range("M62").Value = -1624
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
# Row 65
$ws.Range("H65").Value = 2732.7778
$ws.Range("I65").Value = 2248
$ws.Range("K65").Value = 11240
$ws.Range("M65").Value = -8120
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
# Row 86
$ws.Range("H86").Value = 2827.04
$ws.Range("J86").Value = 3153.6
$ws.Range("L86").Value = 3153.6
$ws.Range("N86").Value = -5399.6
# Row 89
$ws.Range("H89").Value = 2827.04
$ws.Range("J89").Value = 3153.6
$ws.Range("L89").Value = 15768
$ws.Range("N89").Value = -27000
# Row 98
$ws.Range("H98").Value = 1592.04
$ws.Range("I98").Value = 749.65
$ws.Range("K98").Value = 749.65
$ws.Range("M98").Value = 748.35
# Row 106
$ws.Range("H106").Value = 1400
$ws.Range("I106").Value = 850
$ws.Range("K106").Value = 850
$ws.Range("M106").Value = -219
# Row 112
$ws.Range("H112").Value = 1520.6666
$ws.Range("J112").Value = 2025
$ws.Range("L112").Value = 6075
$ws.Range("N112").Value = -8291
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
# Row 122
$ws.Range("H122").Value = 1592.04
$ws.Range("I122").Value = 749.65
$ws.Range("K122").Value = 2248.95
$ws.Range("M122").Value = 201.0500000000002
# Row 137
$ws.Range("H137").Value = 1500.3462
$ws.Range("I137").Value = 1117.6666
$ws.Range("J137").Value = 1702.9412
$ws.Range("K137").Value = 3352.9998
$ws.Range("L137").Value = 5108.8236
$ws.Range("M137").Value = -802.9998000000001
$ws.Range("N137").Value = -10208.8236
# Row 138
$ws.Range("H138").Value = 2848.8628
$ws.Range("J138").Value = 2251.8918
$ws.Range("L138").Value = 6755.6754
$ws.Range("N138").Value = -17035.6754
# Clear cells that no longer have a cached value
$ws.Range("N64").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("M113").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2794
$ws.Range("I2").Value = 2800.842
$ws.Range("J2").Value = 2768
$ws.Range("K2").Value = 2800.842
$ws.Range("L2").Value = 2768
$ws.Range("M2").Value = -2687.842
$ws.Range("N2").Value = -2994
# Row 110
$ws.Range("H110").Value = 1579.9166
$ws.Range("I110").Value = 1483.7142
$ws.Range("K110").Value = 1483.7142
$ws.Range("M110").Value = 561.2858000000001
# Row 116
$ws.Range("H116").Value = 2794
$ws.Range("I116").Value = 2800.842
$ws.Range("J116").Value = 2768
$ws.Range("K116").Value = 2800.842
$ws.Range("L116").Value = 2768
$ws.Range("M116").Value = -506.8420000000001
$ws.Range("N116").Value = -7356
# Row 132
$ws.Range("H132").Value = 2333.9565
$ws.Range("I132").Value = 2085
$ws.Range("J132").Value = 3516.5
$ws.Range("K132").Value = 6255
$ws.Range("L132").Value = 10549.5
$ws.Range("M132").Value = -3725
$ws.Range("N132").Value = -15609.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2794
$ws.Range("I3").Value = 2800.842
$ws.Range("J3").Value = 2768
$ws.Range("K3").Value = 2800.842
$ws.Range("L3").Value = 2768
$ws.Range("M3").Value = -2686.842
$ws.Range("N3").Value = -2996
# Row 20
$ws.Range("H20").Value = 1042.2941
$ws.Range("I20").Value = 1055.6428
$ws.Range("K20").Value = 1055.6428
$ws.Range("M20").Value = -808.6428000000001
# Row 94
$ws.Range("H94").Value = 3450.923
$ws.Range("I94").Value = 942.1818
$ws.Range("K94").Value = 942.1818
$ws.Range("M94").Value = -491.1818
# Row 134
$ws.Range("H134").Value = 2699.0454
$ws.Range("I134").Value = 2206.25
$ws.Range("K134").Value = 6618.75
$ws.Range("M134").Value = -4083.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 2654.6
$ws.Range("I2").Value = 791.3333
$ws.Range("J2").Value = 5449.5
$ws.Range("K2").Value = 791.3333
$ws.Range("L2").Value = 5449.5
$ws.Range("M2").Value = -678.3333
$ws.Range("N2").Value = -5675.5
# Row 31
$ws.Range("H31").Value = 1981.3462
$ws.Range("J31").Value = 2381.2856
$ws.Range("L31").Value = 2381.2856
$ws.Range("N31").Value = -2971.2856
# Row 34
$ws.Range("H34").Value = 1981.3462
$ws.Range("J34").Value = 2381.2856
$ws.Range("L34").Value = 2381.2856
$ws.Range("N34").Value = -2785.2856
# Row 62
$ws.Range("H62").Value = 4239.9287
$ws.Range("I62").Value = 4986.6665
$ws.Range("J62").Value = 2895.8
$ws.Range("K62").Value = 4986.6665
$ws.Range("L62").Value = 2895.8
$ws.Range("M62").Value = -4362.6665
$ws.Range("N62").Value = -4143.8
# Row 65
$ws.Range("H65").Value = 4239.9287
$ws.Range("I65").Value = 4986.6665
$ws.Range("J65").Value = 2895.8
$ws.Range("K65").Value = 24933.3325
$ws.Range("L65").Value = 14479
$ws.Range("M65").Value = -21813.3325
$ws.Range("N65").Value = -20719
# Row 134
$ws.Range("H134").Value = 3331.1155
$ws.Range("I134").Value = 2870.7837
$ws.Range("K134").Value = 8612.3511
$ws.Range("M134").Value = -6077.3511

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 2799.5
$ws.Range("J34").Value = 6700
$ws.Range("L34").Value = 20100
$ws.Range("N34").Value = -20268
# Row 39
$ws.Range("H39").Value = 8305.429
$ws.Range("J39").Value = 8850.462
$ws.Range("L39").Value = 26551.386
$ws.Range("N39").Value = -27139.386
# Row 45
$ws.Range("H45").Value = 5166.6665
$ws.Range("J45").Value = 6750
$ws.Range("L45").Value = 20250
$ws.Range("N45").Value = -21314
# Row 55
$ws.Range("H55").Value = 2882.8333
$ws.Range("J55").Value = 3449.25
$ws.Range("L55").Value = 10347.75
$ws.Range("N55").Value = -10701.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 4022.5
$ws.Range("I9").Value = 4022.5
$ws.Range("K9").Value = 4022.5
$ws.Range("M9").Value = -3852.5
# Row 70
$ws.Range("H70").Value = 7070.615
$ws.Range("I70").Value = 4916
$ws.Range("J70").Value = 10518
$ws.Range("K70").Value = 4916
$ws.Range("L70").Value = 10518
$ws.Range("M70").Value = -4646
$ws.Range("N70").Value = -11058
# Row 73
$ws.Range("H73").Value = 7070.615
$ws.Range("I73").Value = 4916
$ws.Range("J73").Value = 10518
$ws.Range("K73").Value = 4916
$ws.Range("L73").Value = 10518
$ws.Range("M73").Value = -3980
$ws.Range("N73").Value = -12390
# Row 102
$ws.Range("H102").Value = 4170.375
$ws.Range("I102").Value = 3894
$ws.Range("K102").Value = 3894
$ws.Range("M102").Value = -2272
# Row 113
$ws.Range("H113").Value = 2843.5
$ws.Range("J113").Value = 3322.6667
$ws.Range("L113").Value = 3322.6667
$ws.Range("N113").Value = -7662.6667
# Row 122
$ws.Range("H122").Value = 1604.6923
$ws.Range("I122").Value = 1465.091
$ws.Range("K122").Value = 4395.272999999999
$ws.Range("M122").Value = -1945.272999999999
# Row 132
$ws.Range("H132").Value = 3530.1924
$ws.Range("I132").Value = 3149.5557
$ws.Range("J132").Value = 4386.625
$ws.Range("K132").Value = 9448.667099999999
$ws.Range("L132").Value = 13159.875
$ws.Range("M132").Value = -6918.667099999999
$ws.Range("N132").Value = -18219.875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 268.4737
$ws.Range("I16").Value = 268.4737
$ws.Range("K16").Value = 268.4737
$ws.Range("M16").Value = -98.47370000000001
# Row 93
$ws.Range("H93").Value = 4463.6
$ws.Range("I93").Value = 5441.6665
$ws.Range("J93").Value = 2996.5
$ws.Range("K93").Value = 5441.6665
$ws.Range("L93").Value = 2996.5
$ws.Range("M93").Value = -4193.6665
$ws.Range("N93").Value = -5492.5
# Row 136
$ws.Range("H136").Value = 2027.9546
$ws.Range("I136").Value = 1193.125
$ws.Range("J136").Value = 4254.1665
$ws.Range("K136").Value = 3579.375
$ws.Range("L136").Value = 12762.4995
$ws.Range("M136").Value = -1029.375
$ws.Range("N136").Value = -17862.4995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2937.0667
$ws.Range("I122").Value = 2492.2964
$ws.Range("J122").Value = 3604.2222
$ws.Range("K122").Value = 7476.889200000001
$ws.Range("L122").Value = 10812.6666
$ws.Range("M122").Value = -5026.889200000001
$ws.Range("N122").Value = -15712.6666
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
# Clear cells that no longer have a cached value
$ws.Range("N138").ClearContents()

